# "cambios a base de datos" - update the "duracion" (column D) values on the
# "NOdos" sheet for a batch of canales whose average-duration label changed.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("NOdos")

# Rows whose duration text becomes "menos 15 min"
$menos15 = @(21, 23, 27, 28, 29, 34, 35)
foreach ($r in $menos15) {
    $ws.Range("D$r").Value = "menos 15 min"
}

# Rows whose duration text becomes "30 min promedio"
$treinta = @(24, 31)
foreach ($r in $treinta) {
    $ws.Range("D$r").Value = "30 min promedio"
}

$wb.Application.CalculateFull()
